$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.833.12'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.034.85'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'227.24"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').Value = "'0.613"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').Value = "'60.26"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.34%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.386"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').Value = "'0.0816"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').Value = "'14.65"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '2.337.38'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').Value = "'21.07"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('D15').Value = "'0.759"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = "'5.20"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').Value = '2.043.67'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '37.775.13'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').Value = "'6.06"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').Value = "'69.77"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = '0.0₃0824'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('D22').Value = "'225.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -2.47%  '
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('D26').Value = "'9.27"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').Value = "'165.53"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').Value = "'0.130"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.90%  '
$ws.Range('D29').Value = "'18.95"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('E30').Value = '  -6.12%  '
$ws.Range('E31').Value = '  +1.35%  '
$ws.Range('E32').Value = '  -2.80%  '
$ws.Range('E33').Value = '  +4.49%  '
$ws.Range('D34').Value = "'4.52"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.80%  '
$ws.Range('E35').Value = '  -2.38%  '
$ws.Range('D36').Value = "'6.37"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.27%  '
$ws.Range('E37').Value = '  -4.99%  '
$ws.Range('D38').Value = "'3.26"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.17%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '1.539.90'
$ws.Range('E40').Value = '  +3.89%  '
$ws.Range('D41').Value = "'0.0217"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('D42').Value = "'16.95"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('D43').Value = "'97.13"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('D45').Value = "'0.0923"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('D47').Value = "'3.95"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.94%  '
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'7.15"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = "'2.96"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').Value = '2.228.44'
$ws.Range('E51').Value = '  -1.11%  '
